$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Custom Field 1" column header in O1
$ws.Range("O1").Value = "Custom Field 1"

# Fill "Test" for the new custom field column for all data rows (2-9)
$ws.Range("O2:O9").Value = "Test"

# Clear the existing Pan values for the rows that no longer have them
$ws.Range("C3").Value = $null
$ws.Range("C4").Value = $null
$ws.Range("C7").Value = $null
$ws.Range("C9").Value = $null

# Update selection to match the resulting workbook state
$ws.Range("O3:O9").Select()
